# Weekly driver report update for 2025-04-20
# Re-sorts the "Good Drivers" table (rows 12-17) by Driver Vintage (desc)
# and refreshes sample counts / vintage dates for several drivers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a date-like value into column E as TEXT (not an actual
# Excel date serial) while keeping the existing "General" number format
# (style index 4, same as the rest of the table) that column E already
# uses. We do this by temporarily formatting the cell as Text so the
# assignment isn't auto-converted to a date, then restoring the original
# number format by copying it (format-only) from column D on the same
# row, which already carries the correct style.
# ---------------------------------------------------------------------
function Set-TextDate($cell, $sameRowFormatCell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $sameRowFormatCell.Copy()
    $cell.PasteSpecial(-4122)
}

# Row 12: Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B12").Value = 445055
$ws.Range("D12").Value = 99.90000000000001
Set-TextDate $ws.Range("E12") $ws.Range("D12") "2024-11-10"

# Row 13: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B13").Value = 77849
$ws.Range("D13").Value = 99.90000000000001
Set-TextDate $ws.Range("E13") $ws.Range("D13") "2021-08-18"

# Row 14: Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B14").Value = 34244
$ws.Range("D14").Value = 100
Set-TextDate $ws.Range("E14") $ws.Range("D14") "2021-04-27"

# Row 15: Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B15").Value = 59673
$ws.Range("D15").Value = 100
Set-TextDate $ws.Range("E15") $ws.Range("D15") "2020-08-05"

# Row 16: Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B16").Value = 113652
$ws.Range("D16").Value = 100
Set-TextDate $ws.Range("E16") $ws.Range("D16") "2020-01-06"

# Row 17: Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1
# (D17/E17 are unchanged: D17 stays 100, E17 stays "2019-12-14")
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B17").Value = 56018
